$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("moving_average")

# Level (C), Forecast (D), Error (E) for a 4-period moving-average forecast.
$ws.Range("C5").Value = 19500

$ws.Range("C6").Value = 20000
$ws.Range("D6").Value = 9500
$ws.Range("E6").Value = 9500

$ws.Range("C7").Value = 21250
$ws.Range("D7").Value = 2000
$ws.Range("E7").Value = 2000

$ws.Range("C8").Value = 21250
$ws.Range("D8").Value = -1750
$ws.Range("E8").Value = 1750

$ws.Range("C9").Value = 22250
$ws.Range("D9").Value = -16750
$ws.Range("E9").Value = 16750

$ws.Range("C10").Value = 22750
$ws.Range("D10").Value = 10250
$ws.Range("E10").Value = 10250

$ws.Range("C11").Value = 21500
$ws.Range("D11").Value = 9750
$ws.Range("E11").Value = 9750

$ws.Range("C12").Value = 23750
$ws.Range("D12").Value = -10500
$ws.Range("E12").Value = 10500

$ws.Range("C13").Value = 24500
$ws.Range("D13").Value = -17250
$ws.Range("E13").Value = 17250
